$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: B/C value changes from the old "Complementar..." text to the
# instructor string that currently lives in B13:C13. Copy it so the shared
# string + style are reused exactly as Excel would do it natively.
$ws.Range("B13").Copy($ws.Range("B10"))
$ws.Range("C13").Copy($ws.Range("C10"))

# --- Row 18 also ends up with that same instructor string in B/C (copy
# again from the still-untouched B13/C13 before row 13 gets overwritten).
$ws.Range("B13").Copy($ws.Range("B18"))
$ws.Range("C13").Copy($ws.Range("C18"))

# --- Row 15 B/C gets the "01/01/2012" value that already lives in B8:C8.
# Copying (instead of typing the literal string) avoids Excel coercing the
# text into a real date serial number.
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))

# --- Row 23 B/C gets the requirement text currently in B24:C24 (that row
# is being removed below).
$ws.Range("B24").Copy($ws.Range("B23"))
$ws.Range("C24").Copy($ws.Range("C23"))

# --- Column A labels: starting at row 13, every label shifts up by one
# row (new row N label = old row N+1 label). Walk top-down so each read
# happens before its source row is overwritten.
$ws.Range("A14").Copy($ws.Range("A13"))
$ws.Range("A15").Copy($ws.Range("A14"))
$ws.Range("A16").Copy($ws.Range("A15"))
$ws.Range("A17").Copy($ws.Range("A16"))
$ws.Range("A18").Copy($ws.Range("A17"))
$ws.Range("A19").Copy($ws.Range("A18"))
$ws.Range("A20").Copy($ws.Range("A19"))
$ws.Range("A21").Copy($ws.Range("A20"))
$ws.Range("A22").Copy($ws.Range("A21"))
$ws.Range("A23").Copy($ws.Range("A22"))

# --- Row 13 B/C becomes the brand-new "Semestral" text (overwrite after
# the copies above that used B13/C13 as a source).
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# --- Rows that lose their B/C content entirely.
$ws.Range("B14:C14").Clear()
$ws.Range("B16:C16").Clear()
$ws.Range("B22:C22").Clear()

# --- Row 23 no longer carries a label in column A.
$ws.Range("A23").Clear()

# --- Remove the now-surplus row 24 (content already relocated to row 23),
# shifting everything up and shrinking the used range to A1:C23.
$ws.Rows.Item(24).Delete()

# --- Row heights that changed.
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(22).AutoFit()
$ws.Rows.Item(23).RowHeight = 30
